$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -3.900000000000002
$ws.Range("D2").Value = -4.099999999999998
$ws.Range("C3").Value = 4.099999999999998
$ws.Range("D3").Value = 3.5
$ws.Range("C4").Value = 2.699999999999999
$ws.Range("D4").Value = 2
$ws.Range("C5").Value = 0.3999999999999986
$ws.Range("D5").Value = -0.8999999999999986
$ws.Range("C6").Value = -3.900000000000002
$ws.Range("D6").Value = -1.799999999999997
$ws.Range("C7").Value = 0.8000000000000007
$ws.Range("D7").Value = 2.200000000000003
$ws.Range("C8").Value = 1.099999999999998
$ws.Range("D8").Value = 0.3000000000000007
$ws.Range("C9").Value = -1.199999999999999
$ws.Range("D9").Value = -2.299999999999997
$ws.Range("C10").Value = -0.1999999999999993
$ws.Range("D10").Value = 0.2000000000000028
$ws.Range("C11").Value = -1.699999999999999
$ws.Range("D11").Value = 0.2000000000000028
$ws.Range("D12").Value = 2.100000000000001
$ws.Range("C13").Value = 0.3000000000000007
$ws.Range("D13").Value = 1.800000000000001
$ws.Range("C14").Value = -1.5
$ws.Range("D14").Value = 0.7000000000000028
$ws.Range("C15").Value = 4.099999999999998
$ws.Range("D15").Value = -0.6999999999999993
$ws.Range("D16").Value = -2.199999999999999
$ws.Range("D18").Value = 0.5
$ws.Range("C19").Value = -0.8000000000000007
$ws.Range("D19").Value = -4.099999999999998
$ws.Range("D20").Value = -0.2999999999999972
$ws.Range("C21").Value = -2.300000000000001
$ws.Range("D21").Value = -2
$ws.Range("C22").Value = -2.800000000000001
$ws.Range("D22").Value = -1.299999999999997
$ws.Range("C23").Value = -0.1000000000000014
$ws.Range("D23").Value = -3.899999999999999
$ws.Range("C24").Value = -1.5
$ws.Range("D24").Value = 1.100000000000001
$ws.Range("C25").Value = -2.400000000000002
$ws.Range("D25").Value = -0.7999999999999972
$ws.Range("C26").Value = -0.1999999999999993
$ws.Range("D26").Value = -0.5999999999999979
$ws.Range("C27").Value = -0.1000000000000014
$ws.Range("D27").Value = 0.4000000000000021
$ws.Range("C28").Value = -0.4000000000000021
$ws.Range("D28").Value = 3.100000000000001
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = -1.299999999999997
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 2.200000000000003
$ws.Range("C31").Value = 0.3999999999999986
$ws.Range("D31").Value = 1.100000000000001
$ws.Range("C32").Value = 2.599999999999998
$ws.Range("D32").Value = 0.2000000000000028
$ws.Range("C33").Value = 2.300000000000001
$ws.Range("D33").Value = 3.5